$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest years (2004-2009), which occupied rows 2-7.
# Deleting these rows shifts the remaining data (2010-2020, previously
# rows 8-18) up into rows 2-12.
$ws.Range("A2:A7").EntireRow.Delete() | Out-Null

# Append the new 2021 data row as row 13, copying the formatting from the
# row above (A12) so the new label cell keeps the existing header style.
$ws.Range("A12").Copy($ws.Range("A13")) | Out-Null

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value2 = 14.87
$ws.Range("C13").Value2 = 18.84
$ws.Range("D13").Value2 = 98.04000000000001
$ws.Range("E13").Value2 = 99.38
$ws.Range("F13").Value2 = 11.247
$ws.Range("G13").Value2 = 3.29107301607282
